$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original row 13 ("3682251 - Gabrielle Weber Martins" in B/C, no label in A)
# is removed, shifting every row below it up by one (old row 14 -> new row 13, ...,
# old row 22 -> new row 21).
$ws.Rows(13).Delete()

# On top of the shift, several cells now carry different content than what simply
# shifted into place, so overwrite them explicitly to match the final data.

# Row 10 "Objetivos:" now shows the docentes responsaveis text instead of the
# original objectives paragraph.
$ws.Range("B10:C10").Value = "3682251 - Gabrielle Weber Martins"

# Row 13 "Programa resumido:" now just reads "Semestral".
$ws.Range("B13:C13").Value = "Semestral"

# Row 15 "Programa:" now reads the activation date, as text (not as an auto
# converted date serial) -- enter it via a text formula, then convert the
# formula result to a plain value in place so it keeps the original cell
# style/format instead of Excel re-typing the cell as a date.
$ws.Range("B15:C15").Formula = '="01/01/2018"'
$ws.Range("B15:C15").Copy() | Out-Null
$ws.Range("B15:C15").PasteSpecial(-4163) | Out-Null
$ws.Rows(15).RowHeight = 120

# Row 18 "Metodo:" now shows the docentes responsaveis text again.
$ws.Range("B18:C18").Value = "3682251 - Gabrielle Weber Martins"

# Row 19 "Criterio:" now shows the evaluation-method paragraph.
$ws.Range("B19:C19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

# Row 20 "Norma de recuperação:" now shows the previous "Critério" value.
$ws.Range("B20:C20").Value = "NF≥ 5,0."

# Row 21 "Bibliografia:" now shows the previous "Norma de recuperação" value.
$ws.Range("B21:C21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."

Write-Host "Edits applied"
